# The upstream NATMI script was re-run with new TPM input data. The ligand/receptor
# pair (F13a1 -> Itgb1) is unchanged, but the "ECs" sending-cluster block is gone
# from the recomputed output and every remaining specificity/weight column was
# recalculated against the new (smaller) set of edges. Net effect on the sheet:
# the old 9 data rows (3 senders x 3 targets) shrink to 6 (2 senders x 3 targets),
# and most numeric columns get new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old rows 8:10 (the last "MuSCs" sending-cluster block) are no longer needed once
# rows 2:7 are overwritten below with the new 6-row table, so drop them - this also
# shrinks the sheet's used range/dimension from A1:T10 down to A1:T7.
$ws.Rows("8:10").Delete()

# Row 2: FAPs -> F13a1/Itgb1 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "F13a1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07623033333333334
$ws.Range("H2").Value = 0.228691
$ws.Range("I2").Value = 0.7411556909515168
$ws.Range("J2").Value = 0.7411556909515167
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 5.876560060207779
$ws.Range("R2").Value = 52.88904054187001
$ws.Range("S2").Value = 0.1781602407736781
$ws.Range("T2").Value = 0.1781602407736781

# Row 3: FAPs -> F13a1/Itgb1 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "F13a1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07623033333333334
$ws.Range("H3").Value = 0.228691
$ws.Range("I3").Value = 0.7411556909515168
$ws.Range("J3").Value = 0.7411556909515167
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 7.743480105932446
$ws.Range("R3").Value = 69.69132095339201
$ws.Range("S3").Value = 0.2347598367011896
$ws.Range("T3").Value = 0.2347598367011895

# Row 4: FAPs -> F13a1/Itgb1 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F13a1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07623033333333334
$ws.Range("H4").Value = 0.228691
$ws.Range("I4").Value = 0.7411556909515168
$ws.Range("J4").Value = 0.7411556909515167
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 10.82674949314311
$ws.Range("R4").Value = 97.44074543828802
$ws.Range("S4").Value = 0.3282356134766491
$ws.Range("T4").Value = 0.3282356134766491

# Row 5: MuSCs -> F13a1/Itgb1 -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "F13a1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.026623
$ws.Range("H5").Value = 0.079869
$ws.Range("I5").Value = 0.2588443090484832
$ws.Range("J5").Value = 0.2588443090484832
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 2.052354379703333
$ws.Range("R5").Value = 18.47118941733
$ws.Range("S5").Value = 0.0622214265990043
$ws.Range("T5").Value = 0.0622214265990043

# Row 6: MuSCs -> F13a1/Itgb1 -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F13a1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.026623
$ws.Range("H6").Value = 0.079869
$ws.Range("I6").Value = 0.2588443090484832
$ws.Range("J6").Value = 0.2588443090484832
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("Q6").Value = 2.704365333925333
$ws.Range("R6").Value = 24.339288005328
$ws.Range("S6").Value = 0.08198850587687012
$ws.Range("T6").Value = 0.08198850587687011

# Row 7: MuSCs -> F13a1/Itgb1 -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F13a1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.026623
$ws.Range("H7").Value = 0.079869
$ws.Range("I7").Value = 0.2588443090484832
$ws.Range("J7").Value = 0.2588443090484832
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("Q7").Value = 3.781179212421333
$ws.Range("R7").Value = 34.030612911792
$ws.Range("S7").Value = 0.1146343765726088
$ws.Range("T7").Value = 0.1146343765726088
